# Slide 14, content placeholder ("Inhaltsplatzhalter 2"): drop the
# "Improve the documentation" outlook bullet (and the now-superfluous
# trailing empty paragraph after it), and let PowerPoint recompute the
# Shrink-text-on-overflow autofit instead of keeping the stale cached
# fontScale/lnSpcReduction values.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Re-deriving the autofit drops the cached fontScale/lnSpcReduction
# attributes, turning <a:normAutofit fontScale="92500" lnSpcReduction="10000"/>
# back into a plain <a:normAutofit/>.
$tf.AutoSize = 2

# Rebuild the bullet list without the "Improve the documentation" line
# and without the trailing empty paragraph that followed it.
$newText = "Highlights:`r" + `
  "Custom Styling of the Webapp`r" + `
  "Pagination`r" + `
  "React Webapp optimized for mobile devices`r" + `
  "Outlook:`r" + `
  "Planned to do a sprint planning every 2 weeks for continuous releases`r" + `
  "Improve Performance of description search`r" + `
  "Integrate more genes from other data sources"
$tr.Text = $newText

# Restore the sub-bullet indent level that a plain .Text assignment
# resets back to the top level.
$tr.Paragraphs(2).IndentLevel = 2
$tr.Paragraphs(3).IndentLevel = 2
$tr.Paragraphs(4).IndentLevel = 2
$tr.Paragraphs(6).IndentLevel = 2
$tr.Paragraphs(7).IndentLevel = 2
$tr.Paragraphs(8).IndentLevel = 2
